$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 44
$ws_ALC.Range("H44").Value = 47049.25
$ws_ALC.Range("J44").Value = 47049.25
$ws_ALC.Range("L44").Value = 47049.25
$ws_ALC.Range("N44").Value = -47973.25

# ALC row 64
$ws_ALC.Range("H64").Value = 2505.5557
$ws_ALC.Range("I64").Value = 2331.25
$ws_ALC.Range("K64").Value = 2331.25
$ws_ALC.Range("M64").Value = -2083.25

# ALC row 67
$ws_ALC.Range("H67").Value = 2505.5557
$ws_ALC.Range("I67").Value = 2331.25
$ws_ALC.Range("K67").Value = 2331.25
$ws_ALC.Range("M67").Value = -1473.25

# ALC row 70
$ws_ALC.Range("H70").Value = 5079.1665
$ws_ALC.Range("I70").Value = 4993
$ws_ALC.Range("J70").Value = 5096.4
$ws_ALC.Range("K70").Value = 14979
$ws_ALC.Range("L70").Value = 15289.2
$ws_ALC.Range("M70").Value = -14709
$ws_ALC.Range("N70").Value = -15829.2

# ALC row 73
$ws_ALC.Range("H73").Value = 5079.1665
$ws_ALC.Range("I73").Value = 4993
$ws_ALC.Range("J73").Value = 5096.4
$ws_ALC.Range("K73").Value = 14979
$ws_ALC.Range("L73").Value = 15289.2
$ws_ALC.Range("M73").Value = -14043
$ws_ALC.Range("N73").Value = -17161.2

# ALC row 106
$ws_ALC.Range("H106").Value = 69859.60000000001
$ws_ALC.Range("I106").Value = 69859.60000000001
$ws_ALC.Range("K106").Value = 69859.60000000001
$ws_ALC.Range("M106").Value = -69228.60000000001

# ALC row 132
$ws_ALC.Range("H132").Value = 2311.353
$ws_ALC.Range("I132").Value = 2304.25
$ws_ALC.Range("J132").Value = 2344.5
$ws_ALC.Range("K132").Value = 6912.75
$ws_ALC.Range("L132").Value = 7033.5
$ws_ALC.Range("M132").Value = -4382.75
$ws_ALC.Range("N132").Value = -12093.5

# ALC row 135
$ws_ALC.Range("H135").Value = 1219.1428
$ws_ALC.Range("I135").Value = 963.3158
$ws_ALC.Range("K135").Value = 8669.842199999999
$ws_ALC.Range("M135").Value = -6134.842199999999

# ALC row 138
$ws_ALC.Range("H138").Value = 2165.8655
$ws_ALC.Range("I138").Value = 1212.5555
$ws_ALC.Range("K138").Value = 3637.6665
$ws_ALC.Range("M138").Value = 1502.3335

# ARM row 2
$ws_ARM.Range("H2").Value = 2203.7036
$ws_ARM.Range("I2").Value = 2311.6843
$ws_ARM.Range("J2").Value = 1947.25
$ws_ARM.Range("K2").Value = 2311.6843
$ws_ARM.Range("L2").Value = 1947.25
$ws_ARM.Range("M2").Value = -2198.6843
$ws_ARM.Range("N2").Value = -2173.25

# ARM row 37
$ws_ARM.Range("H37").Value = 17800.834
$ws_ARM.Range("J37").Value = 21201.25
$ws_ARM.Range("L37").Value = 21201.25
$ws_ARM.Range("N37").Value = -21747.25

# ARM row 45
$ws_ARM.Range("H45").Value = 4021.8235
$ws_ARM.Range("I45").Value = 4176.3125
$ws_ARM.Range("K45").Value = 4176.3125
$ws_ARM.Range("M45").Value = -3799.3125

# ARM row 74
$ws_ARM.Range("H74").Value = 1218.6666
$ws_ARM.Range("I74").Value = 1215.52
$ws_ARM.Range("J74").Value = 1234.4
$ws_ARM.Range("K74").Value = 1215.52
$ws_ARM.Range("L74").Value = 1234.4
$ws_ARM.Range("M74").Value = -341.52
$ws_ARM.Range("N74").Value = -2982.4

# ARM row 77
$ws_ARM.Range("H77").Value = 1218.6666
$ws_ARM.Range("I77").Value = 1215.52
$ws_ARM.Range("J77").Value = 1234.4
$ws_ARM.Range("K77").Value = 6077.6
$ws_ARM.Range("L77").Value = 6172
$ws_ARM.Range("M77").Value = -1709.6
$ws_ARM.Range("N77").Value = -14908

# ARM row 109
$ws_ARM.Range("H109").Value = 69314.164
$ws_ARM.Range("I109").Value = 0
$ws_ARM.Range("K109").Value = 0
$ws_ARM.Range("M109").ClearContents()

# ARM row 116
$ws_ARM.Range("H116").Value = 2203.7036
$ws_ARM.Range("I116").Value = 2311.6843
$ws_ARM.Range("J116").Value = 1947.25
$ws_ARM.Range("K116").Value = 2311.6843
$ws_ARM.Range("L116").Value = 1947.25
$ws_ARM.Range("M116").Value = -17.68429999999989
$ws_ARM.Range("N116").Value = -6535.25

# BSM row 3
$ws_BSM.Range("H3").Value = 2203.7036
$ws_BSM.Range("I3").Value = 2311.6843
$ws_BSM.Range("J3").Value = 1947.25
$ws_BSM.Range("K3").Value = 2311.6843
$ws_BSM.Range("L3").Value = 1947.25
$ws_BSM.Range("M3").Value = -2197.6843
$ws_BSM.Range("N3").Value = -2175.25

# BSM row 80
$ws_BSM.Range("H80").Value = 531.64703
$ws_BSM.Range("I80").Value = 1574
$ws_BSM.Range("J80").Value = 308.2857
$ws_BSM.Range("K80").Value = 1574
$ws_BSM.Range("L80").Value = 308.2857
$ws_BSM.Range("M80").Value = -576
$ws_BSM.Range("N80").Value = -2304.2857

# BSM row 83
$ws_BSM.Range("H83").Value = 531.64703
$ws_BSM.Range("I83").Value = 1574
$ws_BSM.Range("J83").Value = 308.2857
$ws_BSM.Range("K83").Value = 7870
$ws_BSM.Range("L83").Value = 1541.4285
$ws_BSM.Range("M83").Value = -2878
$ws_BSM.Range("N83").Value = -11525.4285

# BSM row 134
$ws_BSM.Range("H134").Value = 2551.9688
$ws_BSM.Range("I134").Value = 2518.1614
$ws_BSM.Range("K134").Value = 7554.4842
$ws_BSM.Range("M134").Value = -5019.4842

# CRP row 6
$ws_CRP.Range("H6").Value = 6994.1665
$ws_CRP.Range("J6").Value = 4988.6665
$ws_CRP.Range("L6").Value = 4988.6665
$ws_CRP.Range("N6").Value = -5214.6665

# CRP row 31
$ws_CRP.Range("H31").Value = 2115.9578
$ws_CRP.Range("J31").Value = 4265.136
$ws_CRP.Range("L31").Value = 4265.136
$ws_CRP.Range("N31").Value = -4855.136

# CRP row 34
$ws_CRP.Range("H34").Value = 2115.9578
$ws_CRP.Range("J34").Value = 4265.136
$ws_CRP.Range("L34").Value = 4265.136
$ws_CRP.Range("N34").Value = -4669.136

# CRP row 58
$ws_CRP.Range("H58").Value = 6532.826
$ws_CRP.Range("I58").Value = 2187.5557
$ws_CRP.Range("J58").Value = 22175.8
$ws_CRP.Range("K58").Value = 2187.5557
$ws_CRP.Range("L58").Value = 22175.8
$ws_CRP.Range("M58").Value = -1984.5557
$ws_CRP.Range("N58").Value = -22581.8

# CRP row 86
$ws_CRP.Range("H86").Value = 9698
$ws_CRP.Range("I86").Value = 3848
$ws_CRP.Range("K86").Value = 3848
$ws_CRP.Range("M86").Value = -2725

# CRP row 89
$ws_CRP.Range("H89").Value = 9698
$ws_CRP.Range("I89").Value = 3848
$ws_CRP.Range("K89").Value = 19240
$ws_CRP.Range("M89").Value = -13624

# CRP row 132
$ws_CRP.Range("H132").Value = 2361.6
$ws_CRP.Range("I132").Value = 2086.2222
$ws_CRP.Range("J132").Value = 3069.7144
$ws_CRP.Range("K132").Value = 6258.6666
$ws_CRP.Range("L132").Value = 9209.143199999999
$ws_CRP.Range("M132").Value = -3728.6666
$ws_CRP.Range("N132").Value = -14269.1432

# CRP row 136
$ws_CRP.Range("H136").Value = 6532.826
$ws_CRP.Range("I136").Value = 2187.5557
$ws_CRP.Range("J136").Value = 22175.8
$ws_CRP.Range("K136").Value = 6562.6671
$ws_CRP.Range("L136").Value = 66527.39999999999
$ws_CRP.Range("M136").Value = -4012.6671
$ws_CRP.Range("N136").Value = -71627.39999999999

# CUL row 2
$ws_CUL.Range("H2").Value = 156.60976
$ws_CUL.Range("I2").Value = 270.66666
$ws_CUL.Range("K2").Value = 1623.99996
$ws_CUL.Range("M2").Value = -1510.99996

# CUL row 5
$ws_CUL.Range("H5").Value = 1037.1111
$ws_CUL.Range("I5").Value = 639.2308
$ws_CUL.Range("J5").Value = 1406.5714
$ws_CUL.Range("K5").Value = 1917.6924
$ws_CUL.Range("L5").Value = 4219.7142
$ws_CUL.Range("M5").Value = -1805.6924
$ws_CUL.Range("N5").Value = -4443.7142

# CUL row 14
$ws_CUL.Range("H14").Value = 2008.409
$ws_CUL.Range("I14").Value = 2008.409
$ws_CUL.Range("K14").Value = 6025.227000000001
$ws_CUL.Range("M14").Value = -5852.227000000001

# CUL row 107
$ws_CUL.Range("H107").Value = 27779378
$ws_CUL.Range("I107").Value = 404.4
$ws_CUL.Range("J107").Value = 47621500
$ws_CUL.Range("K107").Value = 1213.2
$ws_CUL.Range("L107").Value = 142864500
$ws_CUL.Range("M107").Value = 706.8000000000002
$ws_CUL.Range("N107").Value = -142868340

# CUL row 113
$ws_CUL.Range("H113").Value = 27779174
$ws_CUL.Range("J113").Value = 40001720
$ws_CUL.Range("L113").Value = 120005160
$ws_CUL.Range("N113").Value = -120009500

# CUL row 135
$ws_CUL.Range("H135").Value = 1037.1111
$ws_CUL.Range("I135").Value = 639.2308
$ws_CUL.Range("J135").Value = 1406.5714
$ws_CUL.Range("K135").Value = 5753.077200000001
$ws_CUL.Range("L135").Value = 12659.1426
$ws_CUL.Range("M135").Value = -3218.077200000001
$ws_CUL.Range("N135").Value = -17729.1426

# GSM row 46
$ws_GSM.Range("H46").Value = 29091.072
$ws_GSM.Range("I46").Value = 14999.4
$ws_GSM.Range("J46").Value = 36919.777
$ws_GSM.Range("K46").Value = 14999.4
$ws_GSM.Range("L46").Value = 36919.777
$ws_GSM.Range("M46").Value = -14843.4
$ws_GSM.Range("N46").Value = -37231.777

# GSM row 113
$ws_GSM.Range("H113").Value = 31256022
$ws_GSM.Range("I113").Value = 41673636
$ws_GSM.Range("J113").Value = 3186
$ws_GSM.Range("K113").Value = 41673636
$ws_GSM.Range("L113").Value = 3186
$ws_GSM.Range("M113").Value = -41671466
$ws_GSM.Range("N113").Value = -7526

# GSM row 132
$ws_GSM.Range("H132").Value = 1669488.6
$ws_GSM.Range("I132").Value = 1669488.6
$ws_GSM.Range("J132").Value = 0
$ws_GSM.Range("K132").Value = 5008465.800000001
$ws_GSM.Range("L132").Value = 0
$ws_GSM.Range("M132").Value = -5005935.800000001
$ws_GSM.Range("N132").ClearContents()

# LTW row 61
$ws_LTW.Range("H61").Value = 66671480
$ws_LTW.Range("I61").Value = 100006424
$ws_LTW.Range("J61").Value = 1593.4
$ws_LTW.Range("K61").Value = 100006424
$ws_LTW.Range("L61").Value = 1593.4
$ws_LTW.Range("M61").Value = -100006222
$ws_LTW.Range("N61").Value = -1997.4

# LTW row 113
$ws_LTW.Range("H113").Value = 66671480
$ws_LTW.Range("I113").Value = 100006424
$ws_LTW.Range("J113").Value = 1593.4
$ws_LTW.Range("K113").Value = 100006424
$ws_LTW.Range("L113").Value = 1593.4
$ws_LTW.Range("M113").Value = -100004254
$ws_LTW.Range("N113").Value = -5933.4

# LTW row 122
$ws_LTW.Range("H122").Value = 3313.9092
$ws_LTW.Range("I122").Value = 3205.359
$ws_LTW.Range("K122").Value = 9616.076999999999
$ws_LTW.Range("M122").Value = -7166.076999999999

# LTW row 132
$ws_LTW.Range("H132").Value = 2335488.5
$ws_LTW.Range("I132").Value = 2677840.2
$ws_LTW.Range("K132").Value = 8033520.600000001
$ws_LTW.Range("M132").Value = -8030990.600000001

# LTW row 136
$ws_LTW.Range("H136").Value = 11498317
$ws_LTW.Range("I136").Value = 14496443
$ws_LTW.Range("K136").Value = 43489329
$ws_LTW.Range("M136").Value = -43486779

# WVR row 5
$ws_WVR.Range("H5").Value = 12017801
$ws_WVR.Range("J5").Value = 15017251
$ws_WVR.Range("L5").Value = 15017251
$ws_WVR.Range("N5").Value = -15017475

# WVR row 41
$ws_WVR.Range("H41").Value = 27858.334
$ws_WVR.Range("J41").Value = 28077.4
$ws_WVR.Range("L41").Value = 28077.4
$ws_WVR.Range("N41").Value = -28857.4

# WVR row 100
$ws_WVR.Range("H100").Value = 40000510
$ws_WVR.Range("I100").Value = 62500500
$ws_WVR.Range("J100").Value = 529.55554
$ws_WVR.Range("K100").Value = 125001000
$ws_WVR.Range("L100").Value = 1059.11108
$ws_WVR.Range("M100").Value = -125000459
$ws_WVR.Range("N100").Value = -2141.11108

# WVR row 109
$ws_WVR.Range("H109").Value = 67501.39999999999
$ws_WVR.Range("J109").Value = 67501.39999999999
$ws_WVR.Range("L109").Value = 67501.39999999999
$ws_WVR.Range("N109").Value = -70275.39999999999

# WVR row 113
$ws_WVR.Range("H113").Value = 1527.8
$ws_WVR.Range("I113").Value = 984.1739
$ws_WVR.Range("J113").Value = 3314
$ws_WVR.Range("K113").Value = 3069.7144
$ws_WVR.Range("L113").Value = 9942
$ws_WVR.Range("M113").Value = -782.5217000000002
$ws_WVR.Range("N113").Value = -14282
